# Actualización automática 2025-06-06 08:00:08
$wb = $excel.ActiveWorkbook

$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasGrupo.Range("C4").Value = 497.66
$wsVentasGrupo.Range("L4").Value = 2594.7
$wsVentasGrupo.Range("N35").Value = 1058.37
$wsVentasGrupo.Range("C53").Value = "1 de 51"
$wsVentasGrupo.Range("N53").Value = "1 de 51"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual.Range("F4").Value = 3092.36
$wsVentaMensual.Range("F35").Value = 1058.37
$wsVentaMensual.Range("F53").Value = 7319.76

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento.Range("D2").Value = 497.66
$wsCumplimiento.Range("E2").Value = 9472.68304517915
$wsCumplimiento.Range("F2").Value = 0.0499140298127082

$wsCumplimiento.Range("D16").Value = 5027.97
$wsCumplimiento.Range("E16").Value = 27713.48
$wsCumplimiento.Range("F16").Value = 0.1535658927750604

$wsCumplimiento.Range("D18").Value = 998.24
$wsCumplimiento.Range("E18").Value = 2201.76
$wsCumplimiento.Range("F18").Value = 0.31195

$wsCumplimiento.Range("D19").Value = 7259.63
$wsCumplimiento.Range("E19").Value = 87187.81064517914
$wsCumplimiento.Range("F19").Value = 0.07686423211056648
